# Updates the betexplorer UAE-league scrape sheet:
#   * 6 row pairs had their match data (home/away names, goals, odds,
#     timestamps, url - columns F:V) swapped between two adjacent rows
#     (the "Indice"/pais/torneio/temporada/data_partida columns A:E are
#     untouched because they didn't move).
#   * 2 brand-new match rows were appended at the bottom (rows 49 & 50),
#     growing the used range from A1:V48 to A1:V50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($r1, $r2, $firstCol, $lastCol) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}

# Swap columns F (6) through V (22) for each affected pair of rows.
Swap-RowRange 4 5 6 22
Swap-RowRange 6 7 6 22
Swap-RowRange 18 19 6 22
Swap-RowRange 25 26 6 22
Swap-RowRange 37 38 6 22
Swap-RowRange 41 42 6 22

# Append two new match rows (49 and 50). Copy row 48's formatting down so
# the new rows inherit the same cell styles (bold/bordered index column,
# date-formatted data_partida column) as every other data row.
$ws.Range("A48:V48").Copy($ws.Range("A49:V49"))
$ws.Range("A48:V48").Copy($ws.Range("A50:V50"))

# Row 49: Bani Yas 1 - 0 Khorfakkan
$ws.Cells.Item(49,1).Value2 = 48
$ws.Cells.Item(49,2).Value2 = "united-arab-emirates"
$ws.Cells.Item(49,3).Value2 = "uae-league"
$ws.Cells.Item(49,4).Value2 = "2023-2024"
$ws.Cells.Item(49,5).Value2 = 45234.57291666666
$ws.Cells.Item(49,6).Value2 = "Bani Yas"
$ws.Cells.Item(49,7).Value2 = 1
$ws.Cells.Item(49,8).Value2 = "Khorfakkan"
$ws.Cells.Item(49,9).Value2 = 0
$ws.Cells.Item(49,10).Value2 = 1.9
$ws.Cells.Item(49,11).Value2 = "30/10/2023 18:42"
$ws.Cells.Item(49,12).Value2 = 2
$ws.Cells.Item(49,13).Value2 = "04/11/2023 13:41"
$ws.Cells.Item(49,14).Value2 = 4
$ws.Cells.Item(49,15).Value2 = "30/10/2023 18:42"
$ws.Cells.Item(49,16).Value2 = 3.86
$ws.Cells.Item(49,17).Value2 = "04/11/2023 13:41"
$ws.Cells.Item(49,18).Value2 = 3.56
$ws.Cells.Item(49,19).Value2 = "30/10/2023 18:42"
$ws.Cells.Item(49,20).Value2 = 3.48
$ws.Cells.Item(49,21).Value2 = "04/11/2023 13:41"
$ws.Cells.Item(49,22).Value2 = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/bani-yas-khorfakkan/d8Vq6ptN/"

# Row 50: Ajman 0 - 1 Al Wahda
$ws.Cells.Item(50,1).Value2 = 49
$ws.Cells.Item(50,2).Value2 = "united-arab-emirates"
$ws.Cells.Item(50,3).Value2 = "uae-league"
$ws.Cells.Item(50,4).Value2 = "2023-2024"
$ws.Cells.Item(50,5).Value2 = 45234.6875
$ws.Cells.Item(50,6).Value2 = "Ajman"
$ws.Cells.Item(50,7).Value2 = 0
$ws.Cells.Item(50,8).Value2 = "Al Wahda"
$ws.Cells.Item(50,9).Value2 = 1
$ws.Cells.Item(50,10).Value2 = 4.65
$ws.Cells.Item(50,11).Value2 = "28/10/2023 18:13"
$ws.Cells.Item(50,12).Value2 = 4.61
$ws.Cells.Item(50,13).Value2 = "04/11/2023 16:21"
$ws.Cells.Item(50,14).Value2 = 4.28
$ws.Cells.Item(50,15).Value2 = "28/10/2023 18:13"
$ws.Cells.Item(50,16).Value2 = 4.63
$ws.Cells.Item(50,17).Value2 = "04/11/2023 16:21"
$ws.Cells.Item(50,18).Value2 = 1.64
$ws.Cells.Item(50,19).Value2 = "28/10/2023 18:13"
$ws.Cells.Item(50,20).Value2 = 1.62
$ws.Cells.Item(50,21).Value2 = "04/11/2023 16:21"
$ws.Cells.Item(50,22).Value2 = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ajman-al-wahda/8na2ORei/"
